{"js": "// Add a \"My blog: <url>\" paragraph (as two separate runs) right after the\n// \"For the theory I used the ... attention is all you need paper\" paragraph,\n// replacing the empty paragraph that currently sits just before the table.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph that introduces the theory/reference, then target the\n// paragraph right after it (the empty one immediately preceding the table).\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text.indexOf(\"attention is all you need paper\") !== -1) {\n    target = paragraphs.items[i + 1];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not locate the empty paragraph before the table.\");\n}\n\n// Build the two runs (\"My blog: \" and the URL) as real OOXML so they stay as\n// two distinct <w:r> elements instead of being coalesced into one run.\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body>\" +\n  \"<w:p>\" +\n  '<w:r><w:t xml:space=\"preserve\">My blog: </w:t></w:r>' +\n  \"<w:r><w:t>https://nammibharani.tech/elementor-569/</w:t></w:r>\" +\n  \"</w:p>\" +\n  \"</w:body></w:document>\" +\n  \"</pkg:xmlData></pkg:part></pkg:package>\";\n\ntarget.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Add a \"My blog: <url>\" paragraph (as two separate runs) right after the\n# \"For the theory I used the ... attention is all you need paper\" paragraph,\n# replacing the empty paragraph that currently sits just before the table.\n\n$d = $word.ActiveDocument\n\n# Locate the theory/reference paragraph, then the empty paragraph right after\n# it (the one immediately preceding the table) is our insertion target.\n$searchRange = $d.Content\n$found = $searchRange.Find.Execute(\"attention is all you need paper\")\n$refPara = $searchRange.Paragraphs(1)\n$targetPara = $refPara.Next()\n$target = $targetPara.Range\n\n# Build the two runs (\"My blog: \" and the URL) as real OOXML via InsertXML so\n# they stay as two distinct <w:r> elements instead of being coalesced into a\n# single run.\n$xml = '<?xml version=\"1.0\" standalone=\"yes\"?><?mso-application progid=\"Word.Document\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t xml:space=\"preserve\">My blog: </w:t></w:r><w:r><w:t>https://nammibharani.tech/elementor-569/</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$target.InsertXML($xml)\n"}
